# Add season-record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new header cells AD1:AF1, formatted like the other
# header cells (copy formatting from A1, which carries the bold / centered /
# thin-border header style, then set the text).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2 through 44).
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 94  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 68  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
